$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# STEP 1: Relocate the footer (signature) block from rows 26:27 down
# to rows 33:34 (the data table grows by 7 rows, 21 -> 28, so the
# trailing "NOMBRE DEL REPRESENTANTE LEGAL / FIRMA ..." block that used
# to sit right under the table now needs to sit 7 rows further down,
# leaving a blank gap in between exactly like the original layout did).
# -----------------------------------------------------------------
$ws.Range("B26:C26").UnMerge()
$ws.Range("B27:C27").UnMerge()
$ws.Range("H26:J26").UnMerge()
$ws.Range("H27:J27").UnMerge()

$ws.Range("B26").Copy($ws.Range("B33"))
$ws.Range("C26").Copy($ws.Range("C33"))
$ws.Range("H26").Copy($ws.Range("H33"))
$ws.Range("I26").Copy($ws.Range("I33"))
$ws.Range("J26").Copy($ws.Range("J33"))

$ws.Range("B27").Copy($ws.Range("B34"))
$ws.Range("C27").Copy($ws.Range("C34"))
$ws.Range("H27").Copy($ws.Range("H34"))
$ws.Range("I27").Copy($ws.Range("I34"))
$ws.Range("J27").Copy($ws.Range("J34"))

$ws.Range("B33:C33").Merge()
$ws.Range("B34:C34").Merge()
$ws.Range("H33:J33").Merge()
$ws.Range("H34:J34").Merge()

$ws.Range("B26").Clear()
$ws.Range("C26").Clear()
$ws.Range("H26").Clear()
$ws.Range("I26").Clear()
$ws.Range("J26").Clear()
$ws.Range("B27").Clear()
$ws.Range("C27").Clear()
$ws.Range("H27").Clear()
$ws.Range("I27").Clear()
$ws.Range("J27").Clear()

# -----------------------------------------------------------------
# STEP 2: Grow the worker table from 6 rows (16:21) to 13 rows
# (16:28) by duplicating existing, already-formatted rows into the
# new row slots. This reuses the workbook's existing cell styles
# (normal data row vs. last/bottom-bordered row) instead of inventing
# new ones.
# -----------------------------------------------------------------
# 2a. Preserve the old "last row" (bottom border) styling by copying
#     it down to row 28, which will be the new last row.
$ws.Range("B21:J21").Copy($ws.Range("B28:J28"))

# 2b. Row 21 is no longer the last row, so restyle it like a normal
#     data row (copy format+content down from row 20 first).
$ws.Range("B20:J20").Copy($ws.Range("B21:J21"))

# 2c. Populate the newly-needed normal rows 22:26 from the normal-row
#     block 16:20.
$ws.Range("B16:J20").Copy($ws.Range("B22:J26"))

# 2d. Row 27 also needs normal-row styling.
$ws.Range("B21:J21").Copy($ws.Range("B27:J27"))

# -----------------------------------------------------------------
# STEP 3: Write the real data.
# Rows 16-21: new worker ALVARO JOSE ROJAS CASTRO (CC 8742539),
#             periods 2507 down to 2502.
# Rows 22-28: existing worker HERNAN DARIO SOSA CARDENAS
#             (CC 1054373197), periods 2507 down to 2501 (one more
#             period, 2507, than before).
# -----------------------------------------------------------------
$alvaroRows = @(16,17,18,19,20,21)
$alvaroPeriods = @("2507","2506","2505","2504","2503","2502")
for ($i = 0; $i -lt $alvaroRows.Length; $i++) {
    $r = $alvaroRows[$i]
    $ws.Cells.Item($r, 3).Value = "8742539"
    $ws.Cells.Item($r, 4).Value = "ALVARO JOSE ROJAS CASTRO"
    $ws.Cells.Item($r, 5).Value = $alvaroPeriods[$i]
    $ws.Cells.Item($r, 6).Value = 212128
    $ws.Cells.Item($r, 7).Value = 5303187
}

$hernanRows = @(22,23,24,25,26,27,28)
$hernanPeriods = @("2507","2506","2505","2504","2503","2502","2501")
for ($i = 0; $i -lt $hernanRows.Length; $i++) {
    $r = $hernanRows[$i]
    $ws.Cells.Item($r, 3).Value = "1054373197"
    $ws.Cells.Item($r, 4).Value = "HERNAN DARIO SOSA CARDENAS"
    $ws.Cells.Item($r, 5).Value = $hernanPeriods[$i]
    $ws.Cells.Item($r, 6).Value = 69600
    $ws.Cells.Item($r, 7).Value = 1740000
}

# -----------------------------------------------------------------
# STEP 4: Update the summary cells at the top of the statement.
# -----------------------------------------------------------------
$ws.Range("E11").Value = 1759968   # VALOR MORA (sum of all periods)
$ws.Range("C13").Value = 2         # Cant. Trabajadores
$ws.Range("F13").Value = 7         # Cant. Periodos

Write-Host "done"
